{"js": "// Replace the \"Add free life...\" TODO bullet with the \"Make collider\n// smaller...\" text (which currently lives in the very next bullet), and\n// drop that now-duplicated following paragraph \u2014 i.e. merge the second\n// bullet up into the first one's slot, keeping the first paragraph's\n// properties (list formatting, rsids, the _GoBack bookmark) intact.\n\nconst OLD_TEXT = \"Add free life after 10000 pts (check to make sure) and add free life sound.\";\nconst NEW_TEXT_MARKER = \"Make collider smaller for\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nlet freeLifeIdx = -1;\nlet colliderIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (freeLifeIdx === -1 && t.indexOf(OLD_TEXT) !== -1) freeLifeIdx = i;\n  if (colliderIdx === -1 && t.indexOf(NEW_TEXT_MARKER) !== -1) colliderIdx = i;\n}\n\nif (freeLifeIdx === -1) {\n  throw new Error(\"Could not find the 'Add free life...' paragraph.\");\n}\nif (colliderIdx === -1) {\n  throw new Error(\"Could not find the 'Make collider smaller...' paragraph.\");\n}\n\n// Grab the full text of the \"Make collider smaller...\" bullet, then splice\n// it into the \"Add free life...\" paragraph's slot (keeping that paragraph's\n// own formatting/bookmark), and remove the now-redundant source paragraph.\nconst colliderText = items[colliderIdx].text;\n\nitems[freeLifeIdx].insertText(colliderText, \"Replace\");\nitems[colliderIdx].delete();\n\nawait context.sync();\n", "ps1": "# Replace the \"Add free life...\" TODO bullet with the \"Make collider\n# smaller...\" text (which currently lives in the very next bullet), and\n# drop that now-duplicated following paragraph - i.e. merge the second\n# bullet up into the first one's slot, keeping the first paragraph's own\n# properties (list formatting, rsids, the _GoBack bookmark) intact.\n\n$d = $word.ActiveDocument\n\n$oldTextMarker = \"Add free life after 10000 pts\"\n$newTextMarker = \"Make collider smaller for\"\n\n$freeLifeIndex = -1\n$colliderIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($freeLifeIndex -eq -1 -and $t.Contains($oldTextMarker)) {\n    $freeLifeIndex = $i\n  }\n  if ($colliderIndex -eq -1 -and $t.Contains($newTextMarker)) {\n    $colliderIndex = $i\n  }\n}\n\nif ($freeLifeIndex -eq -1) {\n  throw \"Could not find the 'Add free life...' paragraph.\"\n}\nif ($colliderIndex -eq -1) {\n  throw \"Could not find the 'Make collider smaller...' paragraph.\"\n}\n\n# Full text of the \"Make collider smaller...\" bullet (Range.Text carries the\n# trailing paragraph mark - strip it before using it as a plain replacement).\n$colliderParaRange = $d.Paragraphs.Item($colliderIndex).Range\n$colliderText = $colliderParaRange.Text.TrimEnd([char]13)\n\n$d.Paragraphs.Item($freeLifeIndex).Range.Text = $colliderText\n$d.Paragraphs.Item($colliderIndex).Range.Delete()\n"}
